{"js": "// Replace each two-digit multiplication equation with its new value.\n// The document has exactly one occurrence of each 'old' string, so\n// a direct search+replace per pair is unambiguous.\nconst replacements = [\n  [\"51\u00d768=3468\", \"50\u00d721=1050\"],\n  [\"24\u00d775=1800\", \"24\u00d770=1680\"],\n  [\"92\u00d742=3864\", \"93\u00d731=2883\"],\n  [\"61\u00d727=1647\", \"89\u00d780=7120\"],\n  [\"17\u00d782=1394\", \"68\u00d767=4556\"],\n  [\"76\u00d745=3420\", \"58\u00d758=3364\"],\n  [\"66\u00d771=4686\", \"45\u00d780=3600\"],\n  [\"46\u00d758=2668\", \"23\u00d721=483\"],\n  [\"64\u00d735=2240\", \"55\u00d736=1980\"],\n  [\"96\u00d783=7968\", \"59\u00d740=2360\"],\n  [\"58\u00d738=2204\", \"71\u00d715=1065\"],\n  [\"91\u00d715=1365\", \"66\u00d739=2574\"],\n  [\"61\u00d737=2257\", \"98\u00d764=6272\"],\n  [\"70\u00d749=3430\", \"31\u00d781=2511\"],\n  [\"58\u00d732=1856\", \"18\u00d726=468\"],\n  [\"85\u00d736=3060\", \"97\u00d768=6596\"],\n  [\"69\u00d713=897\", \"95\u00d763=5985\"],\n  [\"94\u00d754=5076\", \"63\u00d732=2016\"],\n  [\"30\u00d726=780\", \"73\u00d717=1241\"],\n  [\"62\u00d773=4526\", \"86\u00d728=2408\"],\n  [\"65\u00d779=5135\", \"45\u00d739=1755\"],\n  [\"50\u00d768=3400\", \"97\u00d779=7663\"],\n  [\"13\u00d743=559\", \"31\u00d751=1581\"],\n  [\"19\u00d726=494\", \"90\u00d764=5760\"],\n  [\"34\u00d763=2142\", \"34\u00d773=2482\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "# Update two-digit multiplication equations throughout the document.\n# Each 'old' equation string occurs exactly once, so Find/Replace is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"51\u00d768=3468\", \"50\u00d721=1050\"),\n    @(\"24\u00d775=1800\", \"24\u00d770=1680\"),\n    @(\"92\u00d742=3864\", \"93\u00d731=2883\"),\n    @(\"61\u00d727=1647\", \"89\u00d780=7120\"),\n    @(\"17\u00d782=1394\", \"68\u00d767=4556\"),\n    @(\"76\u00d745=3420\", \"58\u00d758=3364\"),\n    @(\"66\u00d771=4686\", \"45\u00d780=3600\"),\n    @(\"46\u00d758=2668\", \"23\u00d721=483\"),\n    @(\"64\u00d735=2240\", \"55\u00d736=1980\"),\n    @(\"96\u00d783=7968\", \"59\u00d740=2360\"),\n    @(\"58\u00d738=2204\", \"71\u00d715=1065\"),\n    @(\"91\u00d715=1365\", \"66\u00d739=2574\"),\n    @(\"61\u00d737=2257\", \"98\u00d764=6272\"),\n    @(\"70\u00d749=3430\", \"31\u00d781=2511\"),\n    @(\"58\u00d732=1856\", \"18\u00d726=468\"),\n    @(\"85\u00d736=3060\", \"97\u00d768=6596\"),\n    @(\"69\u00d713=897\", \"95\u00d763=5985\"),\n    @(\"94\u00d754=5076\", \"63\u00d732=2016\"),\n    @(\"30\u00d726=780\", \"73\u00d717=1241\"),\n    @(\"62\u00d773=4526\", \"86\u00d728=2408\"),\n    @(\"65\u00d779=5135\", \"45\u00d739=1755\"),\n    @(\"50\u00d768=3400\", \"97\u00d779=7663\"),\n    @(\"13\u00d743=559\", \"31\u00d751=1581\"),\n    @(\"19\u00d726=494\", \"90\u00d764=5760\"),\n    @(\"34\u00d763=2142\", \"34\u00d773=2482\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $found = $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $found) {\n        throw \"Could not find text: $oldText\"\n    }\n}\n"}
